# "Pruebas realizadas hasta crear cliente recaudo"
#
# - Corrige el nombre de usuario "sgracia" -> "sgarcia" (Santiago Garcia)
# - Agrega dos nombres nuevos a la fila de Nombres (columnas I y J)
# - Agrega las columnas de email para los usuarios ahenao, slopez, storres y
#   sgarcia (columnas E-H de la fila 6), con sus respectivos hipervinculos
#   mailto:
# - Actualiza la celda seleccionada

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New emails for the existing users (typed left to right).
$ws.Range("E6").Value = "ahenao@edeq.com"
$ws.Range("F6").Value = "slopez@gmail.com"
$ws.Range("G6").Value = "storres@hotmail.com"

# Fix the "Nombre Usuario" typo for Santiago Garcia before filling in his
# email.
$ws.Range("H5").Value = "sgarcia"
$ws.Range("H6").Value = "sgarcia@edeq.com"

# Two more names added to the "Nombres" row.
$ws.Range("I3").Value = "Juan David restrepo"
$ws.Range("J3").Value = "Lina maria duran"

# Turn the new email cells into mailto: hyperlinks, same as the existing
# B6:D6 cells.
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:ahenao@edeq.com")
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:slopez@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G6"), "mailto:storres@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("H6"), "mailto:sgarcia@edeq.com")

# Reapply the original hyperlink cell formatting (Hyperlinks.Add resets it to
# a generic style) so E6:H6 match D6 / B6 exactly.
$ws.Range("D6").Copy()
$ws.Range("E6:F6").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("G6:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection where the author last left it.
$ws.Range("B11").Select() | Out-Null
